$d = $word.ActiveDocument

# 1. Date line: "May 09, 2019" -> "January"
$d.Content.Find.Execute("May 09, 2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "January", 2)

# 2. "Hon. Michael P. Cortuna" -> "Hon. "
$d.Content.Find.Execute("Hon. Michael P. Cortuna", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hon. ", 2)

# 3. "Governor" -> "Provincial Director"
$d.Content.Find.Execute("Governor", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Provincial Director", 2)

# 4. "City Of Angeles, Pampanga" -> "Apalit, Pampanga"
$d.Content.Find.Execute("City Of Angeles, Pampanga", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Apalit, Pampanga", 2)

# 5. "Dear Michael P. Cortuna," -> "Dear ,"
$d.Content.Find.Execute("Dear Michael P. Cortuna,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dear ,", 2)

# 6. "FY 2019" -> "FY 2020"
$d.Content.Find.Execute("FY 2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "FY 2020", 2)

# 7. Delete the 7 data rows of the table (keep header row)
$table = $d.Tables.Item(1)
for ($i = $table.Rows.Count; $i -ge 2; $i--) {
    $table.Rows.Item($i).Delete()
}

# 8. "Mark Angelo Maca " -> "President "
$d.Content.Find.Execute("Mark Angelo Maca ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "President ", 2)

# 9. "                                                        Provincial Director" -> "                                                        Michael Cortuna"
$d.Content.Find.Execute("                                                        Provincial Director", $true, $false, $false, $false, $false,
                         $true, 1, $false, "                                                        Michael Cortuna", 2)
